# Exploit for case 10,14,18,22,27
# Fill in the previously-blank "Testcase / Attack script name / Flag / Student ID"
# columns (E:H) for the bug-report rows belonging to case10, case14, case18,
# case22 and case27 — all submitted by student A0127604L.
#
# The assignment-entry order below mirrors how the rows were actually filled
# in (row 26/case10, then row 37/case14, then row 13/case22, then row
# 35/case27, then row 27/case18) so that shared-string de-duplication lines
# up with the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 - case10
$ws.Range("F26").Value = "case10.py"
$ws.Range("G26").Value = "-"

# Row 37 - case14
$ws.Range("F37").Value = "case14.py"
$ws.Range("G37").Value = "EfTj7BxYg2ywfeD"
$ws.Range("H37").Value = "A0127604L"

$ws.Range("H26").Value = "A0127604L"

# Row 13 - case22
$ws.Range("F13").Value = "case22.py"
$ws.Range("G13").Value = "upcYmp7DWrwXF9k"
$ws.Range("H13").Value = "A0127604L"

# Row 35 - case27
$ws.Range("F35").Value = "case27.py"
$ws.Range("G35").Value = "zwPHRtruk8T6S5s"
$ws.Range("H35").Value = "A0127604L"

# Row 27 - case18
$ws.Range("F27").Value = "case18.py"
$ws.Range("G27").Value = "-"
$ws.Range("H27").Value = "A0127604L"

# Testcase numbers (column E)
$ws.Range("E13").Value = 22
$ws.Range("E26").Value = 10
$ws.Range("E27").Value = 18
$ws.Range("E35").Value = 27
$ws.Range("E37").Value = 14
